# Mantel correlogram table update:
#  - number of permutations changed -> narrower "p" column (gridCol width)
#  - recomputed p-values for several distance classes (1000 permutations)
#  - one p-value crossed the significance threshold and is now bolded

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Narrow the 4th (p-value) column: 1090 -> 961 twentieths-of-a-point (dxa)
#    i.e. 54.5pt -> 48.05pt
$t.Columns.Item(4).Width = 48.05

# 2. Plain text-only replacements in the "p" column (row is 1-based, header = row 1)
$plainChanges = @(
    @{ Row = 2;  Old = "<0.001"; New = "0.001" },
    @{ Row = 3;  Old = "0.105";  New = "0.096" },
    @{ Row = 4;  Old = "0.211";  New = "0.192" },
    @{ Row = 5;  Old = "0.422";  New = "0.4"   },
    @{ Row = 6;  Old = "0.845";  New = "0.799" },
    @{ Row = 10; Old = "0.844";  New = "0.783" },
    @{ Row = 11; Old = "0.229";  New = "0.126" },
    @{ Row = 12; Old = "0.185";  New = "0.18"  },
    @{ Row = 13; Old = "0.7";    New = "0.782" },
    @{ Row = 16; Old = "0.008";  New = "0.015" },
    @{ Row = 19; Old = "0.382";  New = "0.336" },
    @{ Row = 23; Old = "0.095";  New = "0.057" },
    @{ Row = 25; Old = "0.026";  New = "0.042" },
    @{ Row = 27; Old = "0.279";  New = "0.3"   },
    @{ Row = 32; Old = "0.346";  New = "0.338" },
    @{ Row = 64; Old = "0.476";  New = "0.521" },
    @{ Row = 65; Old = "0.122";  New = "0.236" },
    @{ Row = 81; Old = "0.123";  New = "0.3"   }
)

# NOTE: Word's Find/Replace engine in this runtime has been observed to
# corrupt unrelated cells when the search text is a textual prefix of a
# value written by an earlier replacement in the same run (e.g. "0.7" is
# a prefix of a previously-written "0.799"). To stay safe, mutate each
# cell's text directly via Range.Text (re-fetching the cell/range fresh
# each time), rather than via Find.Execute.
foreach ($chg in $plainChanges) {
    $cell = $t.Cell($chg.Row, 4)
    $rng = $cell.Range
    $textRng = $d.Range($rng.Start, $rng.End - 1)
    if ($textRng.Text -ne $chg.Old) {
        throw "Unexpected text in row $($chg.Row): got [$($textRng.Text)], expected [$($chg.Old)]"
    }
    $textRng.Text = $chg.New
}

# 3. Row 21 ("9,750" distance class): value changes AND becomes bold (now significant)
$cell = $t.Cell(21, 4)
$rng = $cell.Range
$textRng = $d.Range($rng.Start, $rng.End - 1)
if ($textRng.Text -ne "0.056") {
    throw "Unexpected text in row 21: got [$($textRng.Text)], expected [0.056]"
}
$textRng.Text = "0.02"

$cell = $t.Cell(21, 4)
$rng = $cell.Range
$runRange = $d.Range($rng.Start, $rng.End - 1)
$runRange.Bold = 1
